$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C (rows 2-10) holds a "Förändrad" date that advanced by one day
# (serial 45174 -> 45175, i.e. 2023-09-05 -> 2023-09-06).
for ($row = 2; $row -le 10; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45174) {
        $cell.Value = 45175
    }
}
